$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-17 16:26:44"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-10-17 16:26:21"
$wsZhCn.Range("K2").Value = "2016-10-17 16:27:29"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-10-17 16:26:44"
$wsDeDe.Range("K2").Value = "2016-10-17 16:28:07"
